# UndoRedoNewCommand1StackDiagram.pptx update
#
# The diagram previously illustrated the (now removed) UndoRedoStack /
# UndoableCommand design (AddCommand / DeleteCommand pushing onto a stack).
# It is updated to illustrate the AddressBookCareTaker design instead,
# which keeps a List<ReadOnlyAddressBook> of address book states
# (":AddressBook1" / ":AddressBook0", "prevAddressBook = sN").
#
# This also turns the original pair of "before/after" mini tables into a
# 3-deep "stack" of mini tables per column, so two extra table shapes are
# duplicated from the existing "Table 22" shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) { return $sh }
    }
    return $null
}

function Set-CellText($table, $row, $col, $text) {
    # Table cells here can hold more runs (and, for row 1, a trailing
    # <a:endParaRPr/>) than the new text needs, and a plain
    # `TextRange.Text = "..."` assignment only overwrites the first run,
    # leaving any extra trailing runs in place (e.g. ":AddCommand" would
    # become ":AddressBook1AddCommand"). Forcing a full paragraph replace
    # first (trailing CR) and then setting the real text collapses
    # everything down to a single clean run/paragraph.
    $cell = $table.Cell($row, $col)
    $cr = [char]13
    $cell.Shape.TextFrame.TextRange.Text = "_" + $cr
    $cell.Shape.TextFrame.TextRange.Text = $text
}

# ---------------------------------------------------------------------
# Existing tables: reword the contents from the old UndoableCommand
# vocabulary to the new AddressBookCareTaker vocabulary, and slide them
# down/over a little to make room for the new "stack" entries.
# ---------------------------------------------------------------------

$table22Shape = Get-ShapeByName $s "Table 22"
$table22 = $table22Shape.Table
Set-CellText $table22 1 1 ":AddressBook1"
Set-CellText $table22 2 1 "prevAddressBook = s1"
$table22Shape.Left = 7378562 / 12700.0
$table22Shape.Top = 2322679 / 12700.0
$table22Shape.Width = 2458129 / 12700.0
$table22Shape.Height = 731520 / 12700.0

$table21Shape = Get-ShapeByName $s "Table 21"
$table21 = $table21Shape.Table
Set-CellText $table21 1 1 ":AddressBook0"
Set-CellText $table21 2 1 "prevAddressBook = s0"
$table21Shape.Left = -137322 / 12700.0
$table21Shape.Top = 3099734 / 12700.0
$table21Shape.Width = 2454721 / 12700.0
$table21Shape.Height = 731520 / 12700.0

$table23Shape = Get-ShapeByName $s "Table 23"
$table23 = $table23Shape.Table
Set-CellText $table23 1 1 ":AddressBook0"
Set-CellText $table23 2 1 "prevAddressBook = s0"
$table23Shape.Left = 7384799 / 12700.0
$table23Shape.Top = 3099734 / 12700.0
$table23Shape.Width = 2458129 / 12700.0
$table23Shape.Height = 731520 / 12700.0

# ---------------------------------------------------------------------
# New tables: duplicate "Table 22" (now showing ":AddressBook1" /
# "prevAddressBook = s1") to build out the rest of the two "stacks" of
# mini tables.
# ---------------------------------------------------------------------

# "Table 13": same AddressBook1/s1 entry, placed in the left-hand stack.
$dupRange1 = $table22Shape.Duplicate()
$table13Shape = $dupRange1.Item(1)
$table13Shape.Name = "Table 13"
$table13Shape.Left = -137322 / 12700.0
$table13Shape.Top = 2322679 / 12700.0
$table13Shape.Width = 2458129 / 12700.0
$table13Shape.Height = 731520 / 12700.0

# "Table 14": AddressBook1/s2 entry, placed above "Table 22" in the
# right-hand stack.
$dupRange2 = $table22Shape.Duplicate()
$table14Shape = $dupRange2.Item(1)
$table14Shape.Name = "Table 14"
$table14Shape.Left = 7378561 / 12700.0
$table14Shape.Top = 1537319 / 12700.0
$table14Shape.Width = 2458129 / 12700.0
$table14Shape.Height = 731520 / 12700.0
Set-CellText $table14Shape.Table 2 1 "prevAddressBook = s2"
